# Add block parameter in print_doc_block: append a new instance row
# (Instance3 / Type1 / 0x20000 / 0x10000) below the existing table on
# the "Top" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top")

$ws.Range("A10").Value = "Instance3"
$ws.Range("B10").Value = "Type1"
$ws.Range("C10").Value = "0x20000"
$ws.Range("D10").Value = "0x10000"

$ws.Range("E10").Select() | Out-Null
